# Update "想去人数" (number of people interested) counts for two rows that
# appear on both the "展览" sheet and the aggregated "全部类型" sheet.
#
#   Row 6 (江西·ShiningStaR动漫游戏文化节5th): F6  2129 -> 2131
#   Row 8 (南昌·AP动漫游戏嘉年华):            F8  1290 -> 1293

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($name -eq "展览" -or $name -eq "全部类型") {
        $ws.Range("F6").Value = 2131
        $ws.Range("F8").Value = 1293
    }
}
